# Auto-generated PowerShell Excel COM-interop script
# Applies updated 'F' column (想去人数) values across sheets per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 734
$ws.Range("F3").Value = 14040
$ws.Range("F4").Value = 14040
$ws.Range("F5").Value = 14060
$ws.Range("F7").Value = 1383
$ws.Range("F8").Value = 5817
$ws.Range("F9").Value = 972
$ws.Range("F10").Value = 568
$ws.Range("F12").Value = 375
$ws.Range("F14").Value = 1523
$ws.Range("F15").Value = 422
$ws.Range("F16").Value = 2130
$ws.Range("F17").Value = 1178
$ws.Range("F18").Value = 1790
$ws.Range("F19").Value = 909
$ws.Range("F20").Value = 31
$ws.Range("F21").Value = 2251
$ws.Range("F22").Value = 551
$ws.Range("F23").Value = 791
$ws.Range("F24").Value = 3271
$ws.Range("F26").Value = 305
$ws.Range("F27").Value = 2329
$ws.Range("F28").Value = 75
$ws.Range("F31").Value = 1767
$ws.Range("F32").Value = 1064
$ws.Range("F33").Value = 1345
$ws.Range("F34").Value = 94
$ws.Range("F36").Value = 4670
$ws.Range("F37").Value = 4745
$ws.Range("F38").Value = 293
$ws.Range("F40").Value = 664
$ws.Range("F42").Value = 3263
$ws.Range("F44").Value = 920
$ws.Range("F45").Value = 331
$ws.Range("F46").Value = 87
$ws.Range("F47").Value = 63
$ws.Range("F48").Value = 4407
$ws.Range("F49").Value = 541

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 111
$ws.Range("F7").Value = 81
$ws.Range("F22").Value = 55

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7419
$ws.Range("F3").Value = 215
$ws.Range("F4").Value = 674

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7419
$ws.Range("F3").Value = 734
$ws.Range("F4").Value = 215
$ws.Range("F5").Value = 674
$ws.Range("F7").Value = 14040
$ws.Range("F8").Value = 14040
$ws.Range("F9").Value = 14060
$ws.Range("F11").Value = 1383
$ws.Range("F12").Value = 5817
$ws.Range("F13").Value = 972
$ws.Range("F15").Value = 81
$ws.Range("F16").Value = 422
$ws.Range("F17").Value = 1178
$ws.Range("F18").Value = 1790
$ws.Range("F20").Value = 791
$ws.Range("F21").Value = 3271
$ws.Range("F22").Value = 305
$ws.Range("F23").Value = 75
$ws.Range("F26").Value = 1767
$ws.Range("F32").Value = 1064
$ws.Range("F33").Value = 1345
$ws.Range("F34").Value = 94
$ws.Range("F36").Value = 4670
$ws.Range("F37").Value = 4745
$ws.Range("F38").Value = 293
$ws.Range("F40").Value = 3263
$ws.Range("F42").Value = 920
$ws.Range("F43").Value = 331
$ws.Range("F44").Value = 87
$ws.Range("F45").Value = 63
$ws.Range("F46").Value = 4407
